# Generate Report for Handback
# Fills in the "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" / "Error Detail" columns for the last data
# row (row 8, the fd6dd3d1-... handback) on both the zh-cn and de-de
# report sheets, widens the "Error Detail" column so the new message is
# readable, and links the newly-populated "Latest Target File" cell back
# to the handback markdown file (same target as column A's link).

$wb = $excel.ActiveWorkbook

$targetMdName  = "fd6dd3d1-d867-4f17-9554-25e8289c3af5.md"
$targetMdUrl   = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/784016a023d39340b95aed8bec39cd72dd1c4d96/e2e/fd6dd3d1-d867-4f17-9554-25e8289c3af5.md"
$errorDetail   = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d99cff46fdb49aa976f765a6e1fbdcdb6279c9ba/e2e/fd6dd3d1-d867-4f17-9554-25e8289c3af5.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/784016a023d39340b95aed8bec39cd72dd1c4d96/e2e/fd6dd3d1-d867-4f17-9554-25e8289c3af5.md."

function Apply-HandbackReport($ws, $handbackFile, $handbackDateTime) {
    # Widen column P ("Error Detail") so the long message is readable.
    $ws.Columns.Item(16).ColumnWidth = 39.166666666666664

    # I8 - Latest Target File: link back to the handback .md file.
    $ws.Hyperlinks.Add($ws.Range("I8"), $targetMdUrl, "", "", $targetMdName)
    $ws.Range("I8").Style = "HyperLink"

    # J8 - Latest Handback File (the xlf that was handed back).
    $ws.Range("J8").Value = $handbackFile

    # K8 - Latest Handback DateTime.
    $ws.Range("K8").Value = $handbackDateTime

    # P8 - Error Detail.
    $ws.Range("P8").Value = $errorDetail
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
Apply-HandbackReport $wsZhCn `
    "fd6dd3d1-d867-4f17-9554-25e8289c3af5.648481f9c967c626ca9b81352dbaa2384698ec4c.zh-cn.xlf" `
    "2016-09-05 12:53:47"

$wsDeDe = $wb.Worksheets.Item("de-de")
Apply-HandbackReport $wsDeDe `
    "fd6dd3d1-d867-4f17-9554-25e8289c3af5.648481f9c967c626ca9b81352dbaa2384698ec4c.de-de.xlf" `
    "2016-09-05 12:53:55"
